# Updated symbol list (crypto price/volume refresh + CoinExToken insertion at row 16)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'303.15"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = "'4.36%"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = "'34.85"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "'12.35%"
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.Value = "'5.155"
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.Value = "'4.72%"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = "'4.84%"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.Value = "'2.377"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = "'7.58%"
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.Value = "'8.007"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = "'3.66%"
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.Value = "'3.949"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = "'5.26%"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.Value = "'0.9293"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = "'1.84%"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.Value = "'0.09882"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "'11.02%"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.Value = "'0.1801"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "'6.73%"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.Value = "'0.08660"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.Value = "'5.06%"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Value = "'0.03315"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = "'6.05%"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = "'0.09893"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "'-0.89%"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.Value = "'0.001496"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "'-0.35%"
$cell.Style = "Normal"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$cell = $ws.Range("D16")
$cell.Value = "'0.04577"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = "'0.71%"
$cell.Style = "Normal"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$cell = $ws.Range("D17")
$cell.Value = "'0.005766"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "'-1.65%"
$cell.Style = "Normal"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$cell = $ws.Range("D18")
$cell.Value = "'3.461"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "'-1.00%"
$cell.Style = "Normal"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$cell = $ws.Range("D19")
$cell.Value = "'2.168"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "'3.95%"
$cell.Style = "Normal"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$cell = $ws.Range("D20")
$cell.Value = "'0.3367"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "'1.22%"
$cell.Style = "Normal"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$cell = $ws.Range("D21")
$cell.Value = "'0.1333"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = "'2.74%"
$cell.Style = "Normal"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$cell = $ws.Range("D22")
$cell.Value = "'4.348"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = "'9.41%"
$cell.Style = "Normal"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$cell = $ws.Range("D23")
$cell.Value = "'0.2300"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = "'5.03%"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.Value = "'0.001218"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = "'0.25%"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.Value = "'0.004460"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "'-2.67%"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.Value = "'0.0001301"
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = "'-0.13%"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "'-0.25%"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.Value = "'0.01788"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = "'12.62%"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.Value = "'0.04792"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = "'7.26%"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.Value = "'0.007747"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = "'5.50%"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.Value = "'0.1410"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = "'6.47%"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = "'0.007166"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "'-25.58%"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.Value = "'0.002100"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "'-9.27%"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = "'0.009186"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = "'11.64%"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.Value = "'0.00006124"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "'0.46%"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.Value = "'0.00000000751"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = "'-0.13%"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = "'39.02%"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.Value = "'0.002002"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = "'-0.13%"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.Value = "'0.00002102"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = "'-0.13%"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.Value = "'0.0002001"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = "'-0.13%"
$cell.Style = "Normal"
